$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117 (shifts existing rows 117..222 down to 118..223),
# mirroring the most-recent weekly price observation being prepended to
# this "Ajo" (garlic) price history sheet.
$ws.Rows.Item(117).Insert()

$row = 117
$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value = [DateTime]::FromOADate(44669)
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 100112003
$ws.Cells.Item($row, 7).Value = "Ajo"
$ws.Cells.Item($row, 8).Value = "Chino"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 100
$ws.Cells.Item($row, 11).Value = 20000
$ws.Cells.Item($row, 12).Value = 21000
$ws.Cells.Item($row, 13).Value = 20500
$ws.Cells.Item($row, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item($row, 15).Value = "China"
$ws.Cells.Item($row, 16).Value = 2050
$ws.Cells.Item($row, 17).Value = 10
$ws.Cells.Item($row, 18).Value = "Hortaliza"
